# ---------------------------------------------------------------------------
# TMCAISummitTemplate.docx -- "Add files via upload" commit
#
# The canonical-XML diff for this commit is dominated by artifacts of a
# file having been re-saved/re-uploaded through a different pipeline
# (SharePoint/OneDrive document-library ingestion):
#   * a handful of *unused* namespace declarations (xmlns:oel, xmlns:w16du)
#     and their mc:Ignorable tokens disappear from document.xml /
#     endnotes.xml / footnotes.xml / numbering.xml -- no element in the
#     document actually used those prefixes, so this is pure boilerplate
#     churn from a different Word build's serializer, not a content edit;
#   * the w16cid:durableId="..." attributes on every <w:num> in
#     numbering.xml disappear -- again a serializer-generation artifact,
#     not something exposed on the Word object model;
#   * six new customXml/item*.xml (+ itemProps*.xml) parts appear -- these
#     are SharePoint content-type/metadata parts stamped on upload, not
#     content a document author adds from inside Word;
# none of which correspond to an action a user/macro takes through the
# Word UI/object model -- there is no Find/Replace, style, or formatting
# operation that changes them.
#
# The one change in the diff that *is* a genuine, nameable document
# property is in word/styles.xml: the built-in "Default Paragraph Font"
# character style gains <w:semiHidden/>. In the real Word object model
# this is driven by the Style.Hidden property (the VBA/COM "Hidden"
# property on a Style maps to w:semiHidden in OOXML), so that is the
# operation we drive here.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

$style = $d.Styles("Default Paragraph Font")
try {
    $style.Hidden = $true
} catch {
    # Some hosts don't wire Style.Hidden through to w:semiHidden; ignore.
}
